# Insert a new data row at row 377 (shifts existing rows 377-436 down to 378-437)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("377:377").Insert()

# Populate the newly inserted row with the new record's data
$ws.Range("A377").Value = 9
$ws.Range("B377").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C377").Value = "Metropolitana"
$ws.Range("D377").Value = 45142
$ws.Range("E377").Value = 13
$ws.Range("F377").Value = 100112021
$ws.Range("G377").Value = "Ají"
$ws.Range("H377").Value = "Americana (o)"
$ws.Range("I377").Value = "Primera"
$ws.Range("J377").Value = 25
$ws.Range("K377").Value = 36000
$ws.Range("L377").Value = 38000
$ws.Range("M377").Value = 36960
$ws.Range("N377").Value = "$/caja 25 kilos"
$ws.Range("O377").Value = "Provincia de Limarí"
$ws.Range("P377").Value = 1478
$ws.Range("Q377").Value = 25
$ws.Range("R377").Value = "Hortaliza"
